$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 1503.4736
$ws.Range("I2").Value = 1708.3636
$ws.Range("J2").Value = 1221.75
$ws.Range("K2").Value = 1708.3636
$ws.Range("L2").Value = 1221.75
$ws.Range("M2").Value = -1595.3636
$ws.Range("N2").Value = -1447.75
# Row 12
$ws.Range("H12").Value = 151.44444
$ws.Range("I12").Value = 163
$ws.Range("J12").Value = 128.33333
$ws.Range("K12").Value = 163
$ws.Range("L12").Value = 128.33333
$ws.Range("M12").Value = 7
$ws.Range("N12").Value = -468.33333
# Row 100
$ws.Range("H100").Value = 7848.95
$ws.Range("I100").Value = 7848.95
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 7848.95
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -7307.95
$ws.Range("N100").ClearContents()
# Row 137
$ws.Range("H137").Value = 27700.584
$ws.Range("I137").Value = 52302.668
$ws.Range("K137").Value = 156908.004
$ws.Range("M137").Value = -154358.004

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Range("H4").Value = 399
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 399
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 399
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -631
# Row 5
$ws.Range("H5").Value = 4592.0435
$ws.Range("I5").Value = 5254.15
$ws.Range("J5").Value = 178
$ws.Range("K5").Value = 5254.15
$ws.Range("L5").Value = 178
$ws.Range("M5").Value = -5142.15
$ws.Range("N5").Value = -402
# Row 32
$ws.Range("H32").Value = 20496.02
$ws.Range("I32").Value = 23049.457
$ws.Range("K32").Value = 23049.457
$ws.Range("M32").Value = -22762.457
# Row 88
$ws.Range("H88").Value = 8782.857
$ws.Range("I88").Value = 1927.6666
$ws.Range("J88").Value = 13924.25
$ws.Range("K88").Value = 1927.6666
$ws.Range("L88").Value = 13924.25
$ws.Range("M88").Value = -1521.6666
$ws.Range("N88").Value = -14736.25
# Row 91
$ws.Range("H91").Value = 8782.857
$ws.Range("I91").Value = 1927.6666
$ws.Range("J91").Value = 13924.25
$ws.Range("K91").Value = 1927.6666
$ws.Range("L91").Value = 13924.25
$ws.Range("M91").Value = -523.6666
$ws.Range("N91").Value = -16732.25
# Row 132
$ws.Range("H132").Value = 1101.8387
$ws.Range("I132").Value = 1009.5926
$ws.Range("J132").Value = 1724.5
$ws.Range("K132").Value = 3028.7778
$ws.Range("L132").Value = 5173.5
$ws.Range("M132").Value = -498.7777999999998
$ws.Range("N132").Value = -10233.5

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 4592.0435
$ws.Range("I4").Value = 5254.15
$ws.Range("J4").Value = 178
$ws.Range("K4").Value = 5254.15
$ws.Range("L4").Value = 178
$ws.Range("M4").Value = -5139.15
$ws.Range("N4").Value = -408
# Row 22
$ws.Range("H22").Value = 300
$ws.Range("I22").Value = 300
$ws.Range("K22").Value = 300
$ws.Range("M22").Value = -127

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 176.42857
$ws.Range("I7").Value = 90
$ws.Range("J7").Value = 291.66666
$ws.Range("K7").Value = 90
$ws.Range("L7").Value = 291.66666
$ws.Range("M7").Value = 23
$ws.Range("N7").Value = -517.66666
# Row 31
$ws.Range("H31").Value = 5000971.5
$ws.Range("J31").Value = 1296.3334
$ws.Range("L31").Value = 1296.3334
$ws.Range("N31").Value = -1886.3334
# Row 34
$ws.Range("H34").Value = 5000971.5
$ws.Range("J34").Value = 1296.3334
$ws.Range("L34").Value = 1296.3334
$ws.Range("N34").Value = -1700.3334
# Row 94
$ws.Range("H94").Value = 4799.4
$ws.Range("I94").Value = 4494.5
$ws.Range("K94").Value = 4494.5
$ws.Range("M94").Value = -4043.5
# Row 107
$ws.Range("H107").Value = 1683.1666
$ws.Range("I107").Value = 1524.75
$ws.Range("K107").Value = 1524.75
$ws.Range("M107").Value = 395.25
# Row 122
$ws.Range("H122").Value = 1848.0625
$ws.Range("I122").Value = 1852.3334
$ws.Range("J122").Value = 1842.5714
$ws.Range("K122").Value = 5557.0002
$ws.Range("L122").Value = 5527.7142
$ws.Range("M122").Value = -3107.0002
$ws.Range("N122").Value = -10427.7142
# Row 132
$ws.Range("H132").Value = 202442.4
$ws.Range("I132").Value = 500606
$ws.Range("J132").Value = 3666.6667
$ws.Range("K132").Value = 1501818
$ws.Range("L132").Value = 11000.0001
$ws.Range("M132").Value = -1499288
$ws.Range("N132").Value = -16060.0001

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 176910300
$ws.Range("I4").Value = 127948330
$ws.Range("K4").Value = 383844990
$ws.Range("M4").Value = -383844878
# Row 23
$ws.Range("H23").Value = 393.30768
$ws.Range("I23").Value = 169.16667
$ws.Range("J23").Value = 585.4286
$ws.Range("K23").Value = 507.50001
$ws.Range("L23").Value = 1756.2858
$ws.Range("M23").Value = -272.50001
$ws.Range("N23").Value = -2226.2858
# Row 24
$ws.Range("H24").Value = 1749.5
$ws.Range("J24").Value = 1749.5
$ws.Range("L24").Value = 5248.5
$ws.Range("N24").Value = -5708.5
# Row 64
$ws.Range("H64").Value = 11428.286
$ws.Range("J64").Value = 5000
$ws.Range("L64").Value = 15000
$ws.Range("N64").Value = -15540
# Row 67
$ws.Range("H67").Value = 11428.286
$ws.Range("J67").Value = 5000
$ws.Range("L67").Value = 15000
$ws.Range("N67").Value = -16872
# Row 68
$ws.Range("H68").Value = 4491.0713
$ws.Range("J68").Value = 4491.0713
$ws.Range("L68").Value = 13473.2139
$ws.Range("N68").Value = -15095.2139
# Row 70
$ws.Range("H70").Value = 4944.4375
$ws.Range("I70").Value = 4555.5
$ws.Range("K70").Value = 13666.5
$ws.Range("M70").Value = -13351.5
# Row 71
$ws.Range("H71").Value = 4491.0713
$ws.Range("J71").Value = 4491.0713
$ws.Range("L71").Value = 40419.64169999999
$ws.Range("N71").Value = -48531.64169999999
# Row 73
$ws.Range("H73").Value = 4944.4375
$ws.Range("I73").Value = 4555.5
$ws.Range("K73").Value = 13666.5
$ws.Range("M73").Value = -12574.5
# Row 103
$ws.Range("H103").Value = 807.5
$ws.Range("I103").Value = 708.3333
$ws.Range("K103").Value = 2124.9999
$ws.Range("M103").Value = -1245.9999

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
# Row 52
$ws.Range("H52").Value = 36666
$ws.Range("I52").Value = 36666
$ws.Range("K52").Value = 36666
$ws.Range("M52").Value = -36407
# Row 102
$ws.Range("H102").Value = 2859.875
$ws.Range("I102").Value = 2366.077
$ws.Range("K102").Value = 2366.077
$ws.Range("M102").Value = -744.0770000000002
# Row 104
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
# Row 132
$ws.Range("H132").Value = 2477.121
$ws.Range("I132").Value = 2190
$ws.Range("J132").Value = 3137.5
$ws.Range("K132").Value = 6570
$ws.Range("L132").Value = 9412.5
$ws.Range("M132").Value = -4040
$ws.Range("N132").Value = -14472.5

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 2544997.2
$ws.Range("J2").Value = 2544997.2
$ws.Range("L2").Value = 2544997.2
$ws.Range("N2").Value = -2545221.2
# Row 16
$ws.Range("H16").Value = 2058.25
$ws.Range("I16").Value = 1886
$ws.Range("J16").Value = 2575
$ws.Range("K16").Value = 1886
$ws.Range("L16").Value = 2575
$ws.Range("M16").Value = -1716
$ws.Range("N16").Value = -2915
# Row 51
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
# Row 53
$ws.Range("H53").Value = 70000
$ws.Range("J53").Value = 70000
$ws.Range("L53").Value = 70000
$ws.Range("N53").Value = -71036
# Row 93
$ws.Range("H93").Value = 1241.2273
$ws.Range("I93").Value = 868.9474
$ws.Range("K93").Value = 868.9474
$ws.Range("M93").Value = 379.0526
# Row 104
$ws.Range("H104").Value = 0
$ws.Range("I104").Value = 0
$ws.Range("K104").Value = 0
$ws.Range("M104").ClearContents()
# Row 136
$ws.Range("H136").Value = 3386.1667
$ws.Range("I136").Value = 2519.1333
$ws.Range("K136").Value = 7557.3999
$ws.Range("M136").Value = -5007.3999

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 99
$ws.Range("H99").Value = 65000
$ws.Range("J99").Value = 65000
$ws.Range("L99").Value = 65000
$ws.Range("N99").Value = -70990
# Row 122
$ws.Range("H122").Value = 78065.69
$ws.Range("I122").Value = 104519.79
$ws.Range("K122").Value = 313559.37
$ws.Range("M122").Value = -311109.37
# Row 126
$ws.Range("H126").Value = 4283.1665
$ws.Range("I126").Value = 3173.375
$ws.Range("K126").Value = 9520.125
$ws.Range("M126").Value = -7050.125
# Row 132
$ws.Range("H132").Value = 44308.824
$ws.Range("I132").Value = 56107.69
$ws.Range("K132").Value = 168323.07
$ws.Range("M132").Value = -165793.07
